# Update countries & provincias Spain
# - Refresh COVID-19 figures for several countries (new totals, new cases,
#   active cases, recovered, critical cases, deaths-today, deaths).
# - A handful of countries changed relative rank (by total cases) and so
#   swapped rows with their neighbour: Kazajistan/Oman, Fiyi/Curazao,
#   Seychelles/Montserrat, Papua Nueva Guinea/Islas Virgenes Britanicas.
# - Bump the "last updated" timestamp in the title cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp row
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 12:10"

# Belgica
$ws.Range("B22").Value = 58186
$ws.Range("C22").Value = 125
$ws.Range("D22").Value = 15769
$ws.Range("E22").Value = 32964
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 23
$ws.Range("H22").Value = 9453

# Indonesia
$ws.Range("B35").Value = 25773
$ws.Range("C35").Value = 557
$ws.Range("D35").Value = 7015
$ws.Range("E35").Value = 17185
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 53
$ws.Range("H35").Value = 1573

# Austria
$ws.Range("B44").Value = 16685
$ws.Range("C44").Value = 30
$ws.Range("D44").Value = 15520
$ws.Range("E44").Value = 497
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 668

# Row 54/55 swap: Kazajistan <-> Oman (Oman now ranks above Kazajistan)
$ws.Range("A54").Value = "Oman"
$ws.Range("B54").Value = 10423
$ws.Range("C54").Value = 603
$ws.Range("D54").Value = 2396
$ws.Range("E54").Value = 7985
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 42

$ws.Range("A55").Value = "Kazajistan"
$ws.Range("B55").Value = 10382
$ws.Range("C55").Value = 450
$ws.Range("D55").Value = 5057
$ws.Range("E55").Value = 5288
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 37

# Malasia
$ws.Range("B63").Value = 7762
$ws.Range("C63").Value = 30
$ws.Range("D63").Value = 6330
$ws.Range("E63").Value = 1317
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 115

# Finlandia
$ws.Range("B67").Value = 6826
$ws.Range("C67").Value = 50
$ws.Range("D67").Value = 5500
$ws.Range("E67").Value = 1010
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 316

# Maldivas
$ws.Range("B99").Value = 1591
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 230
$ws.Range("E99").Value = 1356
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 5

# Albania
$ws.Range("B111").Value = 1122
$ws.Range("C111").Value = 23
$ws.Range("D111").Value = 857
$ws.Range("E111").Value = 232
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 33

# Tunez
$ws.Range("B113").Value = 1076
$ws.Range("C113").Value = 5
$ws.Range("D113").Value = 950
$ws.Range("E113").Value = 78
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 48

# Row 198/199 swap: Fiyi <-> Curazao (Curazao now ranks above Fiyi)
$ws.Range("A198").Value = "Curazao"
$ws.Range("B198").Value = 18
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 14
$ws.Range("E198").Value = 3
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Fiyi"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 15
$ws.Range("E199").Value = 3
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# Row 210/211 swap: Seychelles <-> Montserrat (Montserrat now ranks above Seychelles)
$ws.Range("A210").Value = "Montserrat"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 10
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# Row 213/214 swap: Papua Nueva Guinea <-> Islas Virgenes Britanicas
# (Islas Virgenes Britanicas now ranks above Papua Nueva Guinea)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
